# Update the line-loading results (pl_mw) table on Sheet1 with the new
# 380 kV case values. Columns B, C, E, F, G, H, K, M, O for rows 2-25 are
# updated; columns A, D, I, J, L, N are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","E","F","G","H","K","M","O")

$data = @(
    @(2,0.2706980058241015,0.0764738982412041,0.1747333793427259,0.4443680307746121,0.3607337375789115,0.5426354195498106,0.2708818014527594,0.210527886617335,1.745519021933845),
    @(3,0.2368566378903552,0.07382616152859356,0.1640378975501591,0.3878228170618172,0.366148654094026,0.5488033975287294,0.2371206502836003,0.1878383745688694,1.769688185866741),
    @(4,0.2160200724008519,0.07220194998424034,0.1576042791535244,0.3531389305169483,0.3697910100583641,0.5528554810583941,0.2162950458188675,0.1739473171681141,1.785748541313573),
    @(5,0.2075150131102248,0.07154051673494166,0.1550158844627987,0.3390132514313251,0.3713550000765906,0.5545733546926428,0.2077848301140364,0.1682968250009438,1.792599807873103),
    @(6,0.2061019258374301,0.07143071516034638,0.1545880923276428,0.336668177824194,0.3716195085395526,0.5548626298008976,0.2063703055892319,0.1673591848023932,1.793755959951774),
    @(7,0.215905426107895,0.07219302776964298,0.1575692363608212,0.3529483938344953,0.3698117800788339,0.5528783791744729,0.2161803689399875,0.1738710711554532,1.785839699163247),
    @(8,0.2590418247684738,0.07556068428037577,0.171017755607032,0.4248636149813478,0.3625347942994424,0.5447071548209479,0.2592612302739496,0.2026961931559939,1.753599057122244),
    @(9,0.3431531266114121,0.08217404180393828,0.1984596506617038,0.5661985755041457,0.350791474204577,0.5307853107632354,0.3429594564453282,0.2595425812895442,1.70007334564886),
    @(10,0.404636468453532,0.08703562964662126,0.2192898401272174,0.6702781546542269,0.3437137060528528,0.5218379958225654,0.403953316482017,0.3015071444378918,1.666680296470702),
    @(11,0.4325347098079249,0.08924729671923615,0.2289149281233946,0.7176906081379002,0.340832487017515,0.5180457397682687,0.4315883729234997,0.3206423924055173,1.652781385490272),
    @(12,0.4430883848327483,0.09008475316957743,0.2325814222930376,0.7356546913071611,0.3397902901737098,0.5166496798771547,0.4420365541559192,0.3278949305839376,1.647704357887392),
    @(13,0.4408159511650922,0.08990439543099171,0.2317908095304659,0.7317853510981394,0.3400125700151904,0.5169485676770407,0.4397870998906228,0.3263326819134846,1.648789498923065),
    @(14,0.4334031865321037,0.08931619623932363,0.2292161370672972,0.7191683204515869,0.3407457646159813,0.5179300836015202,0.4324482875060482,0.3212389346187834,1.652359961391596),
    @(15,0.4288612321489325,0.08895589819276495,0.2276419056944832,0.7114413442032514,0.3412012356973477,0.5185364983011809,0.4279508714529925,0.3181197027628073,1.654571232872286),
    @(16,0.4028117762244392,0.08689108821862135,0.218663840752626,0.6671810134426437,0.343908827534122,0.5220914229116573,0.4021450082508125,0.3002575162334864,1.667614653661815),
    @(17,0.3868127096670264,0.08562437204332696,0.2131944844661859,0.6400460337125793,0.3456566743766274,0.5243434509702709,0.3862849967367765,0.289311198469214,1.675947529538576),
    @(18,0.3776038375491169,0.08489580451804102,0.2100626986275884,0.6244449056556647,0.346693839321226,0.5256649160901858,0.3771522919978736,0.2830194405258766,1.680861939437662),
    @(19,0.3744847508212104,0.08464912861265361,0.2090047353246973,0.619163680173358,0.3470504701625003,0.5261168335465669,0.3740583394435077,0.2808898962028437,1.682546738960582),
    @(20,0.3885165283534207,0.08575921518480811,0.2137752518449432,0.642933953830422,0.3454673154440542,0.5241010113222799,0.3879744074875191,0.2904760109191216,1.675047896501425),
    @(21,0.4355807925996658,0.08948896661148353,0.2299717906703123,0.7228739723491628,0.3405290801318372,0.5176407032335106,0.4346043303048361,0.3227349177245316,1.651306173452355),
    @(22,0.466276880103635,0.09192622904549808,0.2406836422463243,0.7751780083420101,0.337586525823383,0.5136515921806648,0.4649825589042393,0.3438554459800613,1.636875034319758),
    @(23,0.4498997735099408,0.09062547056808512,0.2349548839650595,0.7472568307830727,0.3391308966240132,0.5157593228132242,0.4487782113525896,0.3325796231424647,1.644477742599207),
    @(24,0.3877462655084685,0.08569825357876937,0.2135126476211653,0.6416283278902171,0.3455528240048551,0.5242105349818047,0.3872106702347651,0.2899493947649603,1.675454235222801),
    @(25,0.3204522763018076,0.08038425704698682,0.1909196774685711,0.5279251897347166,0.3536968416154025,0.5343265486155957,0.320402977929632,0.2441293916629803,1.713513143224333)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $row[$i + 1]
        $ws.Range("$col$r").Value = $val
    }
}
